$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.630.26'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '3.520.88'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.25'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.87'
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("D8").Value = '3.518.92'
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.19'
$ws.Range("E11").Value = '  +7.48%  '
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.40'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("E14").Value = '  -0.93%  '
$ws.Range("D15").Value = '4.096.60'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '611.88'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.531.55'
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '70.710.13'
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.60'
$ws.Range("E21").Value = '  +1.71%  '
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").Value = '  -9.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.36'
$ws.Range("E24").Value = '  +3.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.61'
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("E26").Value = '  -3.40%  '
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.57'
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.11'
$ws.Range("E29").Value = '  +2.90%  '
$ws.Range("E30").Value = '  -2.15%  '
$ws.Range("E31").Value = '  -2.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.07'
$ws.Range("E32").Value = '  -4.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '642.29'
$ws.Range("E33").Value = '  +12.85%  '
$ws.Range("E34").Value = '  -4.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.83'
$ws.Range("E35").Value = '  -2.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.57'
$ws.Range("E36").Value = '  +2.06%  '
$ws.Range("E37").Value = '  -1.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.78'
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0477'
$ws.Range("E39").Value = '  +5.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.87'
$ws.Range("E40").Value = '  -0.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("D43").Value = '0.0₃0745'
$ws.Range("E43").Value = '  +5.99%  '
$ws.Range("D44").Value = '3.375.54'
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("E45").Value = '  -5.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.90'
$ws.Range("E46").Value = '  -2.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.20'
$ws.Range("E47").Value = '  -2.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.56'
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.73'
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("E51").Value = '  -0.01%  '
